$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 180. This pushes the existing rows 180-208
# down to 181-209, matching the diff's dimension change (A1:T208 -> A1:T209).
$ws.Rows("180").Insert()

# Populate the newly inserted row 180 with the new weekly record
# (same static "Chirimoya" / "Macroferia Regional de Talca" template as
# the rest of the subset, with updated date, grade, volume and prices).
$ws.Range("A180").Value = 5
$ws.Range("B180").Value = "Macroferia Regional de Talca"
$ws.Range("C180").Value = "Maule"
$ws.Range("D180").Value = 45258
$ws.Range("E180").Value = 7
$ws.Range("F180").Value = "Fruta"
$ws.Range("G180").Value = 100107
$ws.Range("H180").Value = "Otros"
$ws.Range("I180").Value = 100107002
$ws.Range("J180").Value = "Chirimoya"
$ws.Range("K180").Value = "Cultivar IV Región"
$ws.Range("L180").Value = "Primera"
$ws.Range("M180").Value = 300
$ws.Range("N180").Value = 18000
$ws.Range("O180").Value = 18000
$ws.Range("P180").Value = 18000
$ws.Range("Q180").Value = "$/bandeja 10 kilos"
$ws.Range("R180").Value = "Provincia de Limarí"
$ws.Range("S180").Value = 1800
$ws.Range("T180").Value = 10
